# "final commit of upload excel file"
# Corrections made to the contacts sheet: fix a blank first name, a typo'd
# email domain, trim trailing separators from the Hobbies lists, fix a
# first name and a street value, and tidy up the row heights slightly so
# the (now two-line-capable) header/rows render consistently.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Maria Villa / M G Street contact): first name was blank, email had
# a typo, and the hobbies list had a trailing stray entry.
$ws.Range("B2").Value = "rohan"
$ws.Range("J2").Value = "tintu@gmail.com"
$ws.Range("L2").Value = "Reading ,Drawing"

# Row 3 (dfgdb contact): first name and street corrected, trailing comma
# removed from hobbies list.
$ws.Range("B3").Value = "mini"
$ws.Range("H3").Value = "abcd"
$ws.Range("L3").Value = "Reading ,Writing"

# Slight row-height bump across the header and the two data rows.
$ws.Rows("1:3").RowHeight = 19.5
